$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30 (D1): Comment cleared
$ws.Range("A30").Value = ""

# Row 48 (U1): Comment cleared
$ws.Range("A48").Value = ""

# Row 49 (U2): LY62L2568LL-55LLI -> 3.3V comment, SOT-223 -> SOT223-4, C86781 -> C347229
$ws.Range("A49").Value = "3.3V"
$ws.Range("C49").Value = "SOT223-4"
$ws.Range("D49").Value = "C347229"

# Row 50 (U3): Comment cleared, TSOP-56 -> TSOP-56/14x20x0.5
$ws.Range("A50").Value = ""
$ws.Range("C50").Value = "TSOP-56/14x20x0.5"

# Row 51 (U4): Comment cleared, SOIC-28 -> SOIC-28/300mil
$ws.Range("A51").Value = ""
$ws.Range("C51").Value = "SOIC-28/300mil"

# Row 52 (U5): Comment cleared, TSOP-32/8x14+20x0.5 -> TSOP-32/8x14+20x0.5 b, C261875 -> C261876
$ws.Range("A52").Value = ""
$ws.Range("C52").Value = "TSOP-32/8x14+20x0.5 b"
$ws.Range("D52").Value = "C261876"

# Rows 53-56: Comment cleared
$ws.Range("A53").Value = ""
$ws.Range("A54").Value = ""
$ws.Range("A55").Value = ""
$ws.Range("A56").Value = ""

# View change: scroll position and selection
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("G23").Select()
